# "added colors to rows"
#
# - Rows 6 and 11 (A:J) get a solid orange fill (FFFFCC66)
# - Row 14 (A:J) gets a solid red fill (FFDF5E5E)
# - I14 changes from 0 to 1
# - B19 (inside the merged A19:G19 cell) becomes the boolean FALSE
# - Four FLOOR(...,1,1) formulas drop their redundant third argument

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row fill colors -------------------------------------------------
# FF29A3CC -> R=41  G=163 B=204 -> 41 + 163*256 + 204*65536
$blue = 41 + (163 * 256) + (204 * 65536)
# FFFFCC66 -> R=255 G=204 B=102 -> 255 + 204*256 + 102*65536
$orange = 255 + (204 * 256) + (102 * 65536)
# FFDF5E5E -> R=223 G=94 B=94 -> 223 + 94*256 + 94*65536
$red = 223 + (94 * 256) + (94 * 65536)

# The workbook's palette also picked up this blue swatch (unused by any
# cell in the end state) while the rows below were being colored - apply
# it to a scratch cell and clear it straight back off so the style table
# ends up shaped the same way, then paint the real rows.
$scratch = $ws.Range("ZZ500")
$scratch.Interior.Color = $blue
$scratch.Clear()

$ws.Range("A6:J6").Interior.Color = $orange
$ws.Range("A11:J11").Interior.Color = $orange
$ws.Range("A14:J14").Interior.Color = $red

# --- I14: 0 -> 1 -------------------------------------------------------
$ws.Range("I14").Value = 1

# --- B19: shared string " " -> boolean FALSE ---------------------------
# B19 sits inside the merged range A19:G19, so a normal Value assignment
# is swallowed by the merge (only the top-left cell of a merge can hold
# data). Stage the boolean in a scratch cell and use PasteSpecial (values
# only), which writes straight into the target cell even under a merge.
$ws.Range("ZZ1").Value = $false
$ws.Range("ZZ1").Copy()
$ws.Range("B19").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()

# --- Formula cleanups: FLOOR(x,1,1) -> FLOOR(x,1) -----------------------
$ws.Range("B22").Formula = '=FLOOR(F17/8,1)&"."&FLOOR(MOD(F17,8),1)&"."&(MOD(F17,8)-FLOOR(MOD(F17,8),1))*60'
$ws.Range("B23").Formula = '=FLOOR(H19,1)&"."&(H19-FLOOR(H19,1))*8&".0"'
$ws.Range("B24").Formula = '=FLOOR(I19,1)&"."&(I19-FLOOR(I19,1))*8&".0"'
$ws.Range("B27").Formula = '=FLOOR(K27/8,1)&"."&FLOOR(MOD(K27,8),1)&"."&(MOD(K27,8)-FLOOR(MOD(K27,8),1))*60'
